$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.9036206007003784
$ws.Range("B1").Value = 2.722719192504883
$ws.Range("C1").Value = 4.499489307403564
$ws.Range("D1").Value = 2.169961214065552
$ws.Range("E1").Value = 1.28072988986969
